$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), matching style of H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I/J columns with data values for rows 2-70
$ijValues = @{
    2 = @(7, 7)
    3 = @(6, 6)
    4 = @(7, 7)
    5 = @(3, 4)
    6 = @(6, 7)
    7 = @(1, 1)
    8 = @(7, 7)
    9 = @(7, 7)
    10 = @(7, 7)
    11 = @(6, 6)
    12 = @(6, 7)
    13 = @(7, 7)
    14 = @(8, 8)
    15 = @(10, 10)
    16 = @(7, 8)
    17 = @(7, 8)
    18 = @(10, 10)
    19 = @(9, 10)
    20 = @(6, 7)
    21 = @(10, 10)
    22 = @(7, 8)
    23 = @(8, 8)
    24 = @(6, 7)
    25 = @(8, 9)
    26 = @(8, 8)
    27 = @(9, 9)
    28 = @(9, 9)
    29 = @(8, 8)
    30 = @(7, 7)
    31 = @(6, 6)
    32 = @(6, 6)
    33 = @(6, 7)
    34 = @(5, 5)
    35 = @(8, 8)
    36 = @(6, 6)
    37 = @(7, 7)
    38 = @(8, 8)
    39 = @(9, 9)
    40 = @(7, 7)
    41 = @(6, 6)
    42 = @(8, 8)
    43 = @(8, 8)
    44 = @(5, 5)
    45 = @(5, 5)
    46 = @(7, 7)
    47 = @(7, 7)
    48 = @(8, 8)
    49 = @(7, 7)
    50 = @(8, 8)
    51 = @(9, 9)
    52 = @(8, 8)
    53 = @(7, 7)
    54 = @(7, 7)
    55 = @(7, 7)
    56 = @(7, 7)
    57 = @(7, 7)
    58 = @(7, 7)
    59 = @(7, 7)
    60 = @(9, 9)
    61 = @(7, 7)
    62 = @(7, 7)
    63 = @(8, 8)
    64 = @(6, 6)
    65 = @(7, 7)
    66 = @(5, 5)
    67 = @(4, 4)
    68 = @(4, 4)
    69 = @(6, 6)
    70 = @(3, 3)
}

foreach ($r in $ijValues.Keys) {
    $vals = $ijValues[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
